$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = 'B.Rise call本_M07.Far away（遥不可及）.png'
$ws.Range("B11").Value = 'B.Rise，Unit，Team B'
$ws.Range("C11").Value = '手 在玻璃雾气上画心
窗外春雨淅沥低语
情愫很莫名
我看向你
每次相遇
对视让我 忐忑不已 害怕被你发现
又期待被回应
言语 却总在擦肩瞬间失灵
或许只有收货心碎
才懂距离可贵
尝试幻想每种结尾
再卑微的自我陶醉
Far away
I''m on my way
追逐凄美的误会
直到玫瑰全都枯萎
就不会显得狼狈
Far away
I''m on my way
藏住的真心 应该不算浪费
Far away Away Away
你 像是藏在咖啡杯底
没有被搅开的糖粒
是种只属于
我的惊喜
想要靠近却又担心
破坏了这距离
会不会眼前的诗意就会分崩离析
Far away
I''m on my way
追逐凄美的误会
直到玫瑰全都枯萎
就不会显得狼狈
Far away
I''m on my way
藏住的真心 应该不算浪费
Far away Away Away
如果只把你 留在回忆
就不会叹息
像泡沫不破才美丽
Far away
I''m on my way
追逐凄美的误会
直到玫瑰全都枯萎
就不会显得狼狈
Far away
I''m on my way
藏住的真心 应该不算浪费
Far away Away Away
Far away
I''m on my way
追逐凄美的误会
直到玫瑰全都枯萎
就不会显得狼狈
Far away
I''m on my way
藏住的真心 应该不算浪费
Far away Away Away'
$ws.Range("C11").WrapText = $true
$ws.Range("D11").Value = 'Far Away'
$ws.Rows(11).RowHeight = 17.6

$ws.Range("A12").Value = 'B.Rise call本_M08.Pick Up The PHONE.png'
$ws.Range("B12").Value = 'B.Rise，Unit，Team B'
$ws.Range("C12").Value = 'L-O-V-E
Hey girl
莫非连简单问候
也要反复练习才开口
Do not Do not
不能让机会溜走
天时地利紧要的关头
I know I know
Ring ring ring ring
瞬间周围世界无比安静
Ring ring ring ring
心疯狂预警
双眼紧闭屏住呼吸
献上全部勇气
Tell me baby
准备好坚定的郑重回答
Pick up pick up pick up
这次别再装傻
不需要谨慎地试探想法
Pick it pick it pick it up
Pick it pick it pick up the love
Pick up pick up pick up the phone
Pick up pick up pick up pick up the love
Pick up pick up pick up the phone
Pick up pick up pick up pick up the love
Hey girl
等待接通的每一刻
都期待着
调整到最温柔音色
氛围吻合
或许有默契感应
或许吸引力早已经共鸣
可是心意
果然还是想直接传递oh
恋爱宣言该被认真聆听
Ring ring ring ring
我决心已定
告诉自己大胆出击
幻想拥有魔力
Tell me baby
准备好坚定的郑重回答
Pick up pick up pick up
这次别再装傻
不需要谨慎地试探想法
Pick it pick it pick it up
Pick it pick it pick up the love
正在接通 遐想跃动
此刻也许你也相同
最佳坦白机会 不会错过绝对
利落干脆
心跳已做好准备 now
准备好坚定的郑重回答
Pick up pick up pick up
这次别再装傻
不需要谨慎地试探想法
Pick it pick it pick it up
Pick it pick it pick up the love
Pick up pick up pick up the phone
Pick up pick up pick up pick up the love
Pick up pick up pick up the phone
Pick up pick up pick up pick up the love'
$ws.Range("C12").WrapText = $true
$ws.Range("D12").Value = 'Pick Up The PHONE'
$ws.Rows(12).RowHeight = 17.6

$ws.Range("A13").Value = 'B.Rise call本_M09-Got me mad（失控法则）.png'
$ws.Range("B13").Value = 'B.Rise，Unit，Team B'
$ws.Range("C13").Value = '此刻 微醺的 危险的 为何燥热
在这 微妙的 伪装的 未央旖旎夜色
暧昧气氛袭来 错乱心跳节拍
快一点甩开 多余的拘谨神态 Babe
趁现在深呼吸 沉陷在朦胧里
直白说你目的 怎么做都可以
That''ll be crazy 试探语气泄露蛛丝马迹
Ont to the two
逐渐拉开帷幕
预设好的每一步
恰好的温度
l want you know
眼神交错的时候
这一瞬间已被占有 Baby
（Baby）
One to the two
逐渐混入迷雾
预设好了下一步
模糊了尺度
I want to know
指尖触碰的时候
什么都无需理由 Baby
I want you to
Got me got me got me mad
Got me got me got me mad
Got me got me got me mad mad
Uh Got me got me got me mad
Got me got me got me mad
Got me got me got me mad mad
All l wan All l need
What''s you waiting for
Anyone Anytime
What you wanna meet
花瓣凌乱散开 底线已然摊牌
会屏息等待 狂欢的时刻到来 Babe
遇见你已动心 预见你也倾心
臣服式的着迷 不反抗的进击
一缕优雅轻柔 却能惊醒沉睡野兽 oh
One to the two
逐渐拉开帷幕
预设好的每一步
恰好的温度
l want you know
眼神交错的时候
这一瞬间已被占有 Baby
I want you to
Got me got me got me mad
Got me got me got me mad
Got me got me got me mad mad
Uh
Got me got me got me mad
Got me got me got me mad
Got me got me got me mad mad
Uh oh
被点燃的爱
眩目虚幻的色彩
多嚣张绮丽 burn for me burn for me
映入了眼眸
将所有 无法得到的渴求
和所有 早该舍弃的残留
Burn it if you can
Show me if you can
Gonna crazy with me with me
Wanna go high
One to the two
逐渐混入迷雾
预设好了下一步
模糊了尺度
I want to know
指尖触碰的时候
什么都无需理由 Baby
I want you to
Got me got me got me mad
Got me got me got me mad
Got me got me got me mad mad
Uh Got me got me got me mad
Got me got me got me mad
Got me got me got me mad mad'
$ws.Range("C13").WrapText = $true
$ws.Range("D13").Value = 'Got me mad'
$ws.Rows(13).RowHeight = 17.6

$ws.Range("A14").Value = 'B.Rise call本_M10-草莓味的颜色.png'
$ws.Range("B14").Value = 'B.Rise，Unit，Team B'
$ws.Range("C14").Value = '喜欢后的空格 是草莓成熟的颜色
喜欢后的空格 是每次与你不谋而合
这份憧憬心意的表达
用哪种修辞手法
想装作不经意 语气却很刻意
也许还是缺少勇气 对不起
气氛并非预想中的尴尬
准备的台词却只剩乱码
忽然结巴 简单逻辑也变得好复杂
背对背 快调整呼吸
怎么办 不然先藏好焦急
但我回头看见 你微笑的瞬间
犹豫的 执拗的 纠结 不安的心
全都清零
喜欢后的空格
是每次与你不谋而合 wow wow wow
我终于鼓起勇气抬头看
你坚定的视线也绝不是偶然
喜欢后的空格
是温柔的无可奈何 wow wow wow
任性的 鲜明的 别扭的 可爱的
我要把无限空白写满
形容你的答案
这一次你不用太快回答
否则矜持都没办法
好像你没注意 或者也是有意
将我所有慌张包庇 没关系
煎熬就是对迟疑的惩罚
所以我不想拖延一刹那
因为会怕 直觉有万分之一的误差
背对背 快调整呼吸
怎么办 不然先藏好焦急
但我回头看见 你微笑的瞬间
在意的 执着的 想念的 无比期待的
全都来自你
喜欢后的空格
是草莓成熟后的颜色 wow wow wow
或许时机尚早恋爱也未满
但心动已将甜蜜的红晕渲染
喜欢后的空格
是每首你哼过的歌 wow wow wow
悸动的 轻快的 青涩的 唯美的
我想用关于你的灵感
描绘梦中的浪漫
某一天忽然间 想尝试勇敢做哪怕一点改变
思考我们之间 无数可能性的未来每天
不止单纯的喜欢你（每一面的你）
不设范围的主观题（每一天的你）
一直会故意留白那憧憬心意 ahh
喜欢后的空格
是草莓成熟后的颜色
或许时机尚早恋爱也未满
但心动已将甜蜜的红晕渲染
喜欢后的空格
你希望我写什么呢
真挚的 肤浅的 夸张的 会心一笑的
万千含义 不止于你
喜欢后的空格
是每次与你不谋而合 wow wow wow
我终于鼓起勇气抬头看
你坚定的视线也绝不是偶然
喜欢后的空格
是温柔的无可奈何 wow wow wow
任性的 鲜明的 别扭的 可爱的
我要把无限空白写满
形容你的答案
关于你的答案
属于我们的答案
喜欢后的名词 一直都会有你的名字'
$ws.Range("C14").WrapText = $true
$ws.Range("D14").Value = '草莓味的颜色'
$ws.Rows(14).RowHeight = 17.6

$ws.Range("A15").Value = 'B.Rise call本_M11-Final 奔向落日.png'
$ws.Range("B15").Value = 'B.Rise，16人曲，Team B'
$ws.Range("C15").Value = '地平线被暖色晕染
转眼将逝的绚烂
尽全力奔跑向海岸
倔强依然
想放弃其实很简单
但还是会有期盼
所以我
继续向前哪怕每次都为时已晚
机会运气时间总是不太够
追逐的幻想会忽然变海市蜃楼
可是我也不能再回头
落日被海面吞没陷入黑暗之前
我会再一次来到你的身边
l''ll be l''ll be l''ll be l''ll be
约定好的这天
既然说过了再见那就一定会再见
不顾一切去证明迟来的是诺言
奔向你 奔向你
一起奔赴的终点
曾经我对天空歌唱
多少落寞与迷茫
好像只有去到远方
梦才滚烫
我想要像以前那样
一直站在你身旁
让夕阳
照亮那尘封的吉他和梦想
忘掉所有只管继续奔走
无限拉长的影子直到城市尽头
没有再停下的理由
落日被海面吞没陷入黑暗之前
我会再一次来到你的身边
l''ll be l''ll be l''ll be l''ll be
约定好的这天
既然说过了再见那就一定会再见
不顾一切去证明迟来的是诺言
奔向你 奔向你
一起奔赴的终点
耳边拂过的风也为我伴奏
这次我会毫不犹豫坚定说出口
一定会奔向你最后
落日被海面吞没于交替瞬间
会如约而至出现在你眼前
l''ll be l''ll be l''ll be l''ll be
终点亦是原点
因为说好了再见那就一定要再见
执着的信念并不是谎言
想见你 想见你
想和你 想和你
一起奔向明天'
$ws.Range("C15").WrapText = $true
$ws.Range("D15").Value = '奔向落日'
$ws.Rows(15).RowHeight = 17.6

$ws.Range("A16").Value = 'B.Rise call本_M12-星光环绕的孤岛.png'
$ws.Range("B16").Value = 'B.Rise，16人曲，Team B'
$ws.Range("C16").Value = '像曾经做过的梦那样
一瞬间被点亮了幻想
伫立在舞台中央
被星光一直环绕的地方
无论将来会是怎样
黯淡或者更加闪亮
此刻我要大声尽情歌唱
被黑暗局限的视野
却有微光从未熄灭
直到某刻当我察觉
并非处于孤单结界
寂寞失落的边缘
依然紧握渺小心愿
想要被看见的执念
从未改变
在晚风轻拂的夏夜
第一次任性脱掉了凉鞋
浪花未能将心冷却
依旧还在崩腾不绝
一个人对着海面
忘我沉浸起舞翩翩
妄想会有某天
让梦悄然实现
要经历多少忽视沉默
等待才会有结果
或许不应该思虑太多
偶然间就会遇见流星划过
当看见远处浮现光芒 辉映着不停歇的波浪
会将憧憬和希望 变幻成此刻泪花的形状
呐喊声在耳边回荡 仿佛坚持终有反响
我也拥有了腾跃的能量
像曾经做过的梦那样 一瞬间被点亮了幻想
伫立在舞台中央 被星光一直环绕的地方
无论将来会是怎样 黯淡或者更加闪亮
此刻我要大声尽情歌唱
在晚风轻拂的夏夜
第一次任性脱掉了凉鞋
浪花未能将心冷却
依旧还在崩腾不绝
一个人对着海面
忘我沉浸起舞翩翩
妄想会有某天
让梦悄然实现
要经历多少忽视沉默
等待才会有结果
或许不应该思虑太多
偶然间就会遇见流星划过
当看见远处浮现光芒 辉映着不停歇的波浪
会将憧憬和希望 变幻成此刻泪花的形状
呐喊声在耳边回荡 仿佛坚持终有反响
我也拥有了腾跃的能量
像曾经做过的梦那样 一瞬间被点亮了幻想
伫立在舞台中央 被星光一直环绕的地方
无论将来会是怎样 黯淡或者更加闪亮
此刻我要大声尽情歌唱'
$ws.Range("C16").WrapText = $true
$ws.Range("D16").Value = '星光环绕的孤岛'
$ws.Rows(16).RowHeight = 17.6

$ws.Range("A17").Value = 'B.Rise call本_M13-BETTER THAN BEST.png'
$ws.Range("B17").Value = 'B.Rise，16人曲，Team B'
$ws.Range("C17").Value = '当我茫然自己是为何而来的时候
就会回想坚持这每一步的理由
或许已经不会再有退路在我身后
所以即使要抓住荆棘也不会放手
对于梦想应该怎样去理解
更高处看到的天空是否会有区别
即便等我抵达传说中的那个世界
已经找不到那欣欣向荣的一切
心再渺小的力量也不枯竭
哪怕在无尽长夜
面临荒芜但幸好前路并未断绝
任汗水模糊了视线方向也不会改变
以我全部换一次孤注尝试去超越
越过曾经被认定是最高处的山巅
所有的人都说故事已完结
该付出多少代价才能把历史续写
而如今我好像也已经不在意痛觉
以希望的名义依旧奔腾着热血
心再渺小的力量也不枯竭
哪怕在无尽长夜
我会去到更高处看一看那片天
比谁都抢先一步攀登眺望向更远
就算留给我的仅剩转瞬的时间
也要在最后一刻之前让梦得以成全
顶峰停留不应该是我追求的终点
而是要再一次破茧去飞跃的起点
当我终于笃定自己为何而来的时候
无需多余理由尽管继续远走
当我决心不会留遗憾也不会再回头
好像所谓梦想其实我早已拥有'
$ws.Range("C17").WrapText = $true
$ws.Range("D17").Value = 'Better Than Best'
$ws.Rows(17).RowHeight = 17.6

$ws.Range("A18").Value = 'B.Rise call本_M14-门.png'
$ws.Range("B18").Value = 'B.Rise，16人曲，Team B'
$ws.Range("C18").Value = '我 沐浴冰凉的雨点
破开土睁开眼
只为了 和传说中的人并肩
她 却好像听不见
与我渐行渐远
重重关上我心里那扇门
门那边的人 是传说
门那边的人 让我失落
门那边的人 总沉默
可我听说
门那边 有梦野千里
它真的属于我吗
躲在前人荣光里
让我日渐失去勇气
谁能够开启
为何他们说我不配呢
难道我还不够笨
不够炙热
我 寻找雷霆和闪电
经过冰与火淬炼
找到了 同样失落的你并肩
我 不再怕被灼伤
因为那一束光
不见得 比你的眼更闪亮
门这边的人 感动我
门这边的人 懂我失落
门这边的人 不沉默
让我相信
门那边 有梦野千里
它真的属于我吗
门的这边却有你
让我获得新的勇气
谁能够开启
现在还有什么重要呢
我们早就一样笨
一样炙热
从来就没有 一扇门
只有能不能 触及灵魂
当我去跨过 某扇门
不如宣称
我脚下 是梦野千里
它真的属于我们
继续唱哭一座城
瓦解所有桎梏的门
让我们开启
新的纪元现在就启程
我们还是那样笨
笨到沸腾'
$ws.Range("C18").WrapText = $true
$ws.Range("D18").Value = '门'
$ws.Rows(18).RowHeight = 17.6

$ws.Range("A19").Value = 'B.Rise call本_M15.人海一粟.png'
$ws.Range("B19").Value = 'B.Rise，16人曲，Team B'
$ws.Range("C19").Value = '梦深处 是我渴望的全部 或是片荒芜
身后每一步 消融在无边迷雾
心跳的频度 仿佛依旧不甘愿被现实说服
再向前一步 以终章为序幕
不知不觉被淹没在汹涌的人潮里
连前路视线也都变得不清晰
尝试着拼尽全力 不顾一切拨开拥挤
总有天要从这无望的轮回逃离
在漠然的世界里
逐渐失去了锐利 毫不起眼的自己
像无声叹息 汪洋中的水滴
可是我仍想证明存在的意义
我不怕 即便前方是断壁悬崖
再挣扎 也不停下逞强的步伐
泪水啊 比起虚度蒸发我更想在征途中挥洒
选择哪一种活法 请听我的回答 
Ha Ah～
时间似乎被封印在熙攘的人潮里
被束缚在这了无生机的原地
尝试着拼尽全力 哪怕黑暗没有缝隙
也相信能等到偶然的万分之一
在凝固的世界里
麻木同化的定义 随时都会被代替
想自由呼吸 反向的目的地
才是我作为自己存在的意义
我不怕 即便前方是断壁悬崖
再挣扎 也不停下逞强的步伐
泪水啊 比起虚度蒸发我更想在征途中挥洒
选择哪一种活法 请听我的回答
穿越人海我我不再害怕
未知远方有梦归处等候我抵达
不后悔每一次出发
是我的回答
冲出人海如愿一刹那
会看见梦无限大
曾经自我激励那些话
也早已给出了回答
我不怕 即便前方是断壁悬崖
再挣扎 也不停下逞强的步伐
泪水啊 比起虚度蒸发我更想在征途中挥洒
选择哪一种活法 请听我的回答
给自己的回答'
$ws.Range("C19").WrapText = $true
$ws.Range("D19").Value = '人海一粟'
$ws.Rows(19).RowHeight = 17.6

$ws.Range("A20").Value = 'B.Rise call本_M16.弧线.png'
$ws.Range("B20").Value = 'B.Rise，16人曲，Team B'
$ws.Range("C20").Value = 'Wow
Wow
我曾在泥潭中垂钓星空
你也曾在深渊里打捞彩虹
其实原来我们没有 彼此想象中那么不同
都在追逐自己的梦
黄粱一梦 也要勇 敢向前冲
不辜负心脏的跳动
我们会奔向不同的终点
抬头也在看同一片蓝天
不管时间 给什么历练
曾许下 的誓言 我们不会变
会有一天挣脱重力的茧
绘画出星群迁徙的弧线
像蒲公英 乘着风冒险
我们会 飞很远 比永远还远
Wow
Wow
Dadada～
Dadada～
就算要睁开眼在黑暗之中
也一定会有光点镌刻进眼眸(你在我的眼眸)
未来当然磕磕碰碰 正因我们会一 起度过
所以才更值得感动
黄粱一梦 也要勇 敢向前冲
不辜负心脏的跳动
我们会奔向不同的终点
抬头也在看同一片蓝天
不管时间 给什么历练
曾许下 的誓言 我们不会变
会有一天挣脱重力的茧
绘画出星群迁徙的弧线
像蒲公英 乘着风冒险
我们会 飞很远 比永远还远
Wow
Wow
Wow
（还有新的冒险）
Wow
Wow
Wow
我们会奔向不同的终点
抬头也在看同一片蓝天
不管时间 给什么历练
曾许下 的誓言 我们不会变
会有一天挣脱重力的茧
绘画出星群迁徙的弧线
像蒲公英 乘着风冒险
我们会 飞很远 比永远还远
Wow
Wow
Dadada～
Dadada～
Wow
Wow
Dadada～
Dadada～'
$ws.Range("C20").WrapText = $true
$ws.Range("D20").Value = '弧线'
$ws.Rows(20).RowHeight = 17.6

$ws.Columns(1).ColumnWidth = 46.86
$ws.Range("F19").Select()
